# Descarga el Excel en la ubicacion deseada, luego comprueba si debe descargar
# otro Excel del mismo cuit a ingresar pero distinto cuit representado y luego
# comprueba si debe cerrar sesion o no.
#
# This adds a validation block (columns F:H -> "Anterior" / "Posterior" / "Control")
# to rows 1-10 and duplicates the CUIT row (row 2) into a new row 3 with the same
# "CUIT para ingresar" (A) but matching "CUIT representado" (B) so the control
# formulas can detect consecutive duplicate CUITs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row: new columns F1 "Anterior", G1 "Posterior", H1 "Control"
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)            # xlPasteFormats - reuse header style

$ws.Range("F1").Value = "Anterior"
$ws.Range("G1").Value = "Posterior"

# H1 reuses the header fill/alignment but only has a left/right thin border
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Borders.LineStyle = -4142          # xlLineStyleNone - clear the full box border
$ws.Range("H1").Borders.Item(7).LineStyle = 1      # xlEdgeLeft
$ws.Range("H1").Borders.Item(10).LineStyle = 1     # xlEdgeRight
$ws.Range("H1").Value = "Control"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Row 2: drop the old per-row custom formatting, reuse the plain bordered
#    style already used by the rest of the table (same as row 3+).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Re-enter the values (ClearFormats left them untouched, but make sure).
$ws.Range("A2").Value = 20246845612
$ws.Range("B2").Value = 30715561227
$ws.Range("C2").Value = "Emilia2023"
$ws.Range("D2").Formula = '=LEFT(CELL("filename"),FIND("[",CELL("filename"))-1)'
$ws.Range("E2").Formula = '=D2&"Deudas\"'

# D2:E2 get a plain thin box border + the date number format (border first,
# then format, so the engine folds the format code back onto the built-in
# numFmtId 14 instead of minting a custom one).
$ws.Range("D2:E2").Borders.LineStyle = 1
$ws.Range("D2:E2").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 3) Row 3: duplicate of row 2's "CUIT para ingresar" with the CUIT
#    representado equal to it (same-CUIT follow-up download row).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 20246845612
$ws.Range("B3").Value = 20246845612
$ws.Range("C3").Value = "Emilia2023"
$ws.Range("D3").Formula = '=LEFT(CELL("filename"),FIND("[",CELL("filename"))-1)'
$ws.Range("E3").Formula = '=D3&"Deudas\"'

# Reuse D2's exact format (numFmt14 + border, no alignment) for D3:E3.
$ws.Range("D2").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Columns F:H (Anterior/Posterior/Control) for every data row 2-10.
# ---------------------------------------------------------------------------
$ws.Range("F2").Formula = "=IF(A2=A1,1,0)"
$ws.Range("G2").Formula = "=IF(A2=A3,1,0)"
$ws.Range("H2").Formula = "=F2+G2"

$ws.Range("F3").Formula = "=IF(A3=A2,1,0)"
$ws.Range("G3").Formula = "=IF(A3=A4,1,0)"
$ws.Range("H3").Formula = "=F3+G3"

for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = $null
    $ws.Cells.Item($r, 7).Value = $null
    $ws.Cells.Item($r, 8).Value = $null
}

# Apply the plain bordered style to every F:H cell, rows 2-10.
$ws.Range("F2:H10").Borders.LineStyle = 1

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) View bookkeeping to mirror the author's final selection/scroll state.
# ---------------------------------------------------------------------------
$ws.Range("G18").Select()
try { $excel.ActiveWindow.ScrollColumn = 4 } catch {}
